$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "27.377.52"
$ws.Cells.Item(2, 5).Value = "  -2.88%  "

Set-TextValue $ws.Cells.Item(3, 4) "1.742.72"
$ws.Cells.Item(3, 5).Value = "  -3.27%  "

$ws.Cells.Item(4, 5).Value = "  -0.01%  "

Set-TextValue $ws.Cells.Item(5, 4) "321.55"
$ws.Cells.Item(5, 5).Value = "  -4.49%  "

Set-TextValue $ws.Cells.Item(6, 4) "1.000"
$ws.Cells.Item(6, 5).Value = "  +0.02%  "

Set-TextValue $ws.Cells.Item(7, 4) "0.4194"
$ws.Cells.Item(7, 5).Value = "  -8.98%  "

Set-TextValue $ws.Cells.Item(8, 4) "0.3575"
$ws.Cells.Item(8, 5).Value = "  -3.36%  "

Set-TextValue $ws.Cells.Item(9, 4) "45.42"
$ws.Cells.Item(9, 5).Value = "  +0.67%  "

Set-TextValue $ws.Cells.Item(10, 4) "0.07403"
$ws.Cells.Item(10, 5).Value = "  -2.44%  "

Set-TextValue $ws.Cells.Item(11, 4) "1.112"
$ws.Cells.Item(11, 5).Value = "  -3.10%  "

$ws.Cells.Item(12, 5).Value = "  -0.12%  "

Set-TextValue $ws.Cells.Item(13, 4) "21.46"
$ws.Cells.Item(13, 5).Value = "  -3.86%  "

Set-TextValue $ws.Cells.Item(14, 4) "6.071"
$ws.Cells.Item(14, 5).Value = "  -4.13%  "

Set-TextValue $ws.Cells.Item(15, 4) "7.166"
$ws.Cells.Item(15, 5).Value = "  -3.58%  "

Set-TextValue $ws.Cells.Item(16, 4) "1.738.34"
$ws.Cells.Item(16, 5).Value = "  -3.44%  "

Set-TextValue $ws.Cells.Item(17, 4) "0.00001065"
$ws.Cells.Item(17, 5).Value = "  -2.90%  "

Set-TextValue $ws.Cells.Item(18, 4) "87.09"
$ws.Cells.Item(18, 5).Value = "  +6.09%  "

Set-TextValue $ws.Cells.Item(19, 4) "0.06022"
$ws.Cells.Item(19, 5).Value = "  -10.43%  "

$ws.Cells.Item(20, 5).Value = "  +0.05%  "

Set-TextValue $ws.Cells.Item(21, 4) "16.83"
$ws.Cells.Item(21, 5).Value = "  -4.03%  "

Set-TextValue $ws.Cells.Item(22, 4) "6.088"
$ws.Cells.Item(22, 5).Value = "  -5.01%  "

Set-TextValue $ws.Cells.Item(23, 4) "0.5235"
$ws.Cells.Item(23, 5).Value = "  -5.90%  "

Set-TextValue $ws.Cells.Item(24, 4) "27.405.01"
$ws.Cells.Item(24, 5).Value = "  -2.73%  "

Set-TextValue $ws.Cells.Item(25, 4) "11.40"
$ws.Cells.Item(25, 5).Value = "  -4.02%  "

Set-TextValue $ws.Cells.Item(26, 4) "2.344"
$ws.Cells.Item(26, 5).Value = "  -2.64%  "

Set-TextValue $ws.Cells.Item(27, 4) "20.39"
$ws.Cells.Item(27, 5).Value = "  -1.40%  "

Set-TextValue $ws.Cells.Item(28, 4) "153.09"

Set-TextValue $ws.Cells.Item(29, 4) "2.378"
$ws.Cells.Item(29, 5).Value = "  +0.09%  "

Set-TextValue $ws.Cells.Item(30, 4) "1.935.28"
$ws.Cells.Item(30, 5).Value = "  -3.49%  "

Set-TextValue $ws.Cells.Item(31, 4) "125.59"
$ws.Cells.Item(31, 5).Value = "  -5.79%  "

Set-TextValue $ws.Cells.Item(32, 4) "1.176"
$ws.Cells.Item(32, 5).Value = "  -6.14%  "

Set-TextValue $ws.Cells.Item(33, 4) "5.664"
$ws.Cells.Item(33, 5).Value = "  -3.22%  "

Set-TextValue $ws.Cells.Item(34, 4) "0.09131"
$ws.Cells.Item(34, 5).Value = "  -4.59%  "

Set-TextValue $ws.Cells.Item(35, 4) "3.613"
$ws.Cells.Item(35, 5).Value = "  -10.39%  "

Set-TextValue $ws.Cells.Item(36, 4) "12.61"
$ws.Cells.Item(36, 5).Value = "  +4.93%  "

Set-TextValue $ws.Cells.Item(37, 4) "0.02284"
$ws.Cells.Item(37, 5).Value = "  -2.84%  "

Set-TextValue $ws.Cells.Item(38, 4) "0.2130"
$ws.Cells.Item(38, 5).Value = "  -4.15%  "

Set-TextValue $ws.Cells.Item(39, 4) "5.064"
$ws.Cells.Item(39, 5).Value = "  -3.56%  "

Set-TextValue $ws.Cells.Item(40, 4) "0.06041"
$ws.Cells.Item(40, 5).Value = "  -5.05%  "

Set-TextValue $ws.Cells.Item(41, 4) "0.6373"
$ws.Cells.Item(41, 5).Value = "  -3.89%  "

Set-TextValue $ws.Cells.Item(42, 4) "1.190"
$ws.Cells.Item(42, 5).Value = "  -3.53%  "

Set-TextValue $ws.Cells.Item(43, 4) "1.430"
$ws.Cells.Item(43, 5).Value = "  -5.15%  "

Set-TextValue $ws.Cells.Item(44, 4) "0.9996"
$ws.Cells.Item(44, 5).Value = "  -0.01%  "

Set-TextValue $ws.Cells.Item(45, 4) "7.902"
$ws.Cells.Item(45, 5).Value = "  -2.32%  "

Set-TextValue $ws.Cells.Item(46, 4) "13.70"
$ws.Cells.Item(46, 5).Value = "  -3.11%  "

Set-TextValue $ws.Cells.Item(47, 4) "3.700"
$ws.Cells.Item(47, 5).Value = "  -3.50%  "

Set-TextValue $ws.Cells.Item(48, 4) "0.5821"
$ws.Cells.Item(48, 5).Value = "  -4.52%  "

Set-TextValue $ws.Cells.Item(49, 4) "125.27"
$ws.Cells.Item(49, 5).Value = "  -3.64%  "

Set-TextValue $ws.Cells.Item(50, 4) "1.941"
$ws.Cells.Item(50, 5).Value = "  -5.14%  "

Set-TextValue $ws.Cells.Item(51, 4) "0.06820"
$ws.Cells.Item(51, 5).Value = "  -4.52%  "
